# Atualização de bases das ligas, do dia: 03-03-2024 às 00:35
#
# This script:
#  1. Swaps the full row content (columns B:AC) between pairs of rows whose
#     match identity (B = game id) was recorded against the wrong row.
#     Column A (the running index) stays put.
#  2. Fills in results (H/I/J) and refreshed odds (K:AC) for three fixtures
#     that have since been played (rows 209-211).
#  3. Refreshes a handful of odds values for rows 212-216.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows {
    param($RowA, $RowB)
    $rangeA = $ws.Range("B$RowA`:AC$RowA")
    $rangeB = $ws.Range("B$RowB`:AC$RowB")
    $valsA = $rangeA.Value()
    $valsB = $rangeB.Value()
    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

# --- 1. Row swaps (game rows that were attributed to the wrong slot) ---
Swap-Rows 33 34
Swap-Rows 66 67
Swap-Rows 70 71
Swap-Rows 76 77
Swap-Rows 87 88
Swap-Rows 94 95
Swap-Rows 97 98
Swap-Rows 139 140
Swap-Rows 167 168
Swap-Rows 195 196
Swap-Rows 204 205

# --- 2. Rows 209-211: fixtures now played, add results + refreshed odds ---

# Row 209: Moreirense 0 - 0 Rio Ave -> Draw
$ws.Range("H209").Value2 = 0
$ws.Range("I209").Value2 = 0
$ws.Range("J209").Value2 = "D"
$ws.Range("N209").Value2 = 2.45
$ws.Range("P209").Value2 = 3.2
$ws.Range("R209").Value2 = 2.1
$ws.Range("S209").Value2 = 1.775
$ws.Range("U209").Value2 = 1.925
$ws.Range("V209").Value2 = 1.925
$ws.Range("W209").Value2 = -1
$ws.Range("X209").Value2 = 2.1
$ws.Range("Y209").Value2 = -1
$ws.Range("Z209").Value2 = -0.5
$ws.Range("AA209").Value2 = 0.3875
$ws.Range("AB209").Value2 = -1
$ws.Range("AC209").Value2 = 0.925

# Row 210: Estoril 1 - 3 Guimaraes -> Away win
$ws.Range("H210").Value2 = 1
$ws.Range("I210").Value2 = 3
$ws.Range("J210").Value2 = "A"
$ws.Range("N210").Value2 = 2.875
$ws.Range("O210").Value2 = 3.4
$ws.Range("P210").Value2 = 2.45
$ws.Range("Q210").Value2 = 0
$ws.Range("R210").Value2 = 2.1
$ws.Range("S210").Value2 = 1.775
$ws.Range("U210").Value2 = 1.875
$ws.Range("V210").Value2 = 1.975
$ws.Range("W210").Value2 = -1
$ws.Range("X210").Value2 = -1
$ws.Range("Y210").Value2 = 1.45
$ws.Range("Z210").Value2 = -1
$ws.Range("AA210").Value2 = 0.7749999999999999
$ws.Range("AB210").Value2 = 0.875
$ws.Range("AC210").Value2 = -1

# Row 211: Braga 3 - 0 Estrela -> Home win
$ws.Range("H211").Value2 = 3
$ws.Range("I211").Value2 = 0
$ws.Range("J211").Value2 = "H"
$ws.Range("P211").Value2 = 8
$ws.Range("Q211").Value2 = -1.25
$ws.Range("R211").Value2 = 1.825
$ws.Range("S211").Value2 = 2.025
$ws.Range("U211").Value2 = 1.925
$ws.Range("V211").Value2 = 1.925
$ws.Range("W211").Value2 = 0.363
$ws.Range("X211").Value2 = -1
$ws.Range("Y211").Value2 = -1
$ws.Range("Z211").Value2 = 0.825
$ws.Range("AA211").Value2 = -1
$ws.Range("AB211").Value2 = 0
$ws.Range("AC211").Value2 = 0

# --- 3. Rows 212-216: refreshed odds on still-upcoming fixtures ---

# Row 212
$ws.Range("U212").Value2 = 1.825
$ws.Range("V212").Value2 = 2.025

# Row 213
$ws.Range("N213").Value2 = 2.55
$ws.Range("P213").Value2 = 2.875
$ws.Range("R213").Value2 = 1.83
$ws.Range("S213").Value2 = 2.07

# Row 214
$ws.Range("N214").Value2 = 1.142
$ws.Range("O214").Value2 = 8
$ws.Range("Q214").Value2 = -2.25
$ws.Range("R214").Value2 = 1.85
$ws.Range("S214").Value2 = 2.05

# Row 215
$ws.Range("R215").Value2 = 2.08
$ws.Range("S215").Value2 = 1.82
$ws.Range("U215").Value2 = 1.85
$ws.Range("V215").Value2 = 2

# Row 216
$ws.Range("R216").Value2 = 1.99
$ws.Range("S216").Value2 = 1.91
